$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the date as plain text "04.03.2018" (shared string) -> "08.03.2018".
# Prefix with an apostrophe so the text isn't auto-converted into a real date value,
# then restore the "Normal" style so no number-format/style gets attached to the cells
# (matching the original, which has no style on A1:A3).
$ws.Range("A1:A3").Value = "'08.03.2018"
$ws.Range("A1:A3").Style = "Normal"

# Column B holds time-of-day values (numFmtId 21 / h:mm:ss) that shift by roughly 40
# minutes. Assign the exact serial-day fractions so the underlying numeric values match.
$ws.Range("B1").Value = 0.10196759259259258
$ws.Range("B2").Value = 0.10197916666666666
$ws.Range("B3").Value = 0.10199074074074073
